$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2808873333333333
$ws.Range("H2").Value = 0.842662
$ws.Range("I2").Value = 0.5595554696739399
$ws.Range("J2").Value = 0.5595554696739399
$ws.Range("M2").Value = 0.8317113333333332
$ws.Range("N2").Value = 2.495134
$ws.Range("O2").Value = 0.0263454906755698
$ws.Range("P2").Value = 0.0263454906755698
$ws.Range("Q2").Value = 0.2336171785231111
$ws.Range("R2").Value = 2.102554606708
$ws.Range("S2").Value = 0.01474176340875886
$ws.Range("T2").Value = 0.01474176340875887

# Row 3
$ws.Range("G3").Value = 0.2808873333333333
$ws.Range("H3").Value = 0.842662
$ws.Range("I3").Value = 0.5595554696739399
$ws.Range("J3").Value = 0.5595554696739399
$ws.Range("O3").Value = 0.6529848313028861
$ws.Range("P3").Value = 0.6529848313028862
$ws.Range("Q3").Value = 5.790306803768444
$ws.Range("R3").Value = 52.112761233916
$ws.Range("S3").Value = 0.3653812339696448
$ws.Range("T3").Value = 0.3653812339696449

# Row 4
$ws.Range("G4").Value = 0.2808873333333333
$ws.Range("H4").Value = 0.842662
$ws.Range("I4").Value = 0.5595554696739399
$ws.Range("J4").Value = 0.5595554696739399
$ws.Range("M4").Value = 10.12334933333333
$ws.Range("N4").Value = 30.370048
$ws.Range("O4").Value = 0.3206696780215441
$ws.Range("P4").Value = 0.3206696780215441
$ws.Range("Q4").Value = 2.843520598641778
$ws.Range("R4").Value = 25.591685387776
$ws.Range("S4").Value = 0.1794324722955362
$ws.Range("T4").Value = 0.1794324722955362

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2210956666666667
$ws.Range("H5").Value = 0.663287
$ws.Range("I5").Value = 0.4404445303260602
$ws.Range("J5").Value = 0.4404445303260602
$ws.Range("M5").Value = 0.8317113333333332
$ws.Range("N5").Value = 2.495134
$ws.Range("O5").Value = 0.0263454906755698
$ws.Range("P5").Value = 0.0263454906755698
$ws.Range("Q5").Value = 0.1838877717175555
$ws.Range("R5").Value = 1.654989945458
$ws.Range("S5").Value = 0.01160372726681094
$ws.Range("T5").Value = 0.01160372726681094

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2210956666666667
$ws.Range("H6").Value = 0.663287
$ws.Range("I6").Value = 0.4404445303260602
$ws.Range("J6").Value = 0.4404445303260602
$ws.Range("O6").Value = 0.6529848313028861
$ws.Range("P6").Value = 0.6529848313028862
$ws.Range("Q6").Value = 4.557741097796222
$ws.Range("R6").Value = 41.019669880166
$ws.Range("S6").Value = 0.2876035973332414
$ws.Range("T6").Value = 0.2876035973332414

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2210956666666667
$ws.Range("H7").Value = 0.663287
$ws.Range("I7").Value = 0.4404445303260602
$ws.Range("J7").Value = 0.4404445303260602
$ws.Range("M7").Value = 10.12334933333333
$ws.Range("N7").Value = 30.370048
$ws.Range("O7").Value = 0.3206696780215441
$ws.Range("P7").Value = 0.3206696780215441
$ws.Range("Q7").Value = 2.238228669752889
$ws.Range("R7").Value = 20.144058027776
$ws.Range("S7").Value = 0.1412372057260079
$ws.Range("T7").Value = 0.141237205726008

